# Swap the content of the specific cells between row 3 and row 4.
# Columns involved: A, B, D, E, F, G, H, M, Q, R, S
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "M", "Q", "R", "S")

foreach ($col in $cols) {
    $cell3 = $ws.Range($col + "3")
    $cell4 = $ws.Range($col + "4")

    $val3 = $cell3.Value()
    $val4 = $cell4.Value()

    $cell3.Value = $val4
    $cell4.Value = $val3
}

# M3 had a value ("äldre gnagspår") and M4 was empty before the edit.
# After swapping, M4 should hold the text and M3 should become empty.
# Since M4 was empty (no value), setting $cell3.Value = $val4 where
# $val4 is empty/null may not actually clear M3's content depending on
# how the COM shim treats empty values, so explicitly clear it to match
# the diff (M3 cell is removed entirely in the output).
$ws.Range("M3").Value = ""
